# Bugfix: recompute the simulation columns (r0, Aktuell Infizierte,
# Dunkelziffer, Aktuell Immune, Aktuell Verstorbene) for rows 2-66 on
# Sheet1 with corrected values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a true 2-D .NET array (object[,]) rather than a PowerShell array-of-
# arrays: Range.Value only accepts a rectangular SAFEARRAY, same as VBA.
$data = New-Object 'object[,]' 65,5
$data[0,0] = 0
$data[0,1] = 0.01666666666666667
$data[0,2] = 1.666666666666667
$data[0,3] = 0
$data[0,4] = 0
$data[1,0] = 0
$data[1,1] = 0.01666666666666667
$data[1,2] = 1.2
$data[1,3] = 0.003333333333333334
$data[1,4] = 0
$data[2,0] = 0
$data[2,1] = 0.02333333333333333
$data[2,2] = 1.142857142857143
$data[2,3] = 0.003333333333333334
$data[2,4] = 0
$data[3,0] = 0
$data[3,1] = 0.03
$data[3,2] = 1
$data[3,3] = 0.003333333333333334
$data[3,4] = 0
$data[4,0] = 0
$data[4,1] = 0.03333333333333333
$data[4,2] = 1
$data[4,3] = 0.006666666666666667
$data[4,4] = 0
$data[5,0] = 0
$data[5,1] = 0.05666666666666666
$data[5,2] = 1
$data[5,3] = 0.006666666666666667
$data[5,4] = 0
$data[6,0] = 0
$data[6,1] = 0.07333333333333333
$data[6,2] = 1
$data[6,3] = 0.01
$data[6,4] = 0
$data[7,0] = 0
$data[7,1] = 0.1033333333333333
$data[7,2] = 1
$data[7,3] = 0.01
$data[7,4] = 0
$data[8,0] = 0
$data[8,1] = 0.1366666666666667
$data[8,2] = 1
$data[8,3] = 0.01
$data[8,4] = 0.003333333333333334
$data[9,0] = 0.6666666666666666
$data[9,1] = 0.1666666666666667
$data[9,2] = 1
$data[9,3] = 0.01666666666666667
$data[9,4] = 0.003333333333333334
$data[10,0] = 0.75
$data[10,1] = 0.2366666666666667
$data[10,2] = 1
$data[10,3] = 0.02333333333333333
$data[10,4] = 0.003333333333333334
$data[11,0] = 1.153846153846154
$data[11,1] = 0.2633333333333333
$data[11,2] = 1
$data[11,3] = 0.03666666666666667
$data[11,4] = 0.006666666666666667
$data[12,0] = 1.142857142857143
$data[12,1] = 0.3033333333333333
$data[12,2] = 1
$data[12,3] = 0.06
$data[12,4] = 0.01
$data[13,0] = 1.269230769230769
$data[13,1] = 0.3166666666666667
$data[13,2] = 1
$data[13,3] = 0.07666666666666666
$data[13,4] = 0.01
$data[14,0] = 1.25
$data[14,1] = 0.3233333333333333
$data[14,2] = 1
$data[14,3] = 0.09666666666666666
$data[14,4] = 0.01
$data[15,0] = 1.230769230769231
$data[15,1] = 0.3533333333333333
$data[15,2] = 1
$data[15,3] = 0.12
$data[15,4] = 0.01
$data[16,0] = 1.18
$data[16,1] = 0.37
$data[16,2] = 1
$data[16,3] = 0.1533333333333333
$data[16,4] = 0.01333333333333333
$data[17,0] = 1.392857142857143
$data[17,1] = 0.3833333333333334
$data[17,2] = 1
$data[17,3] = 0.17
$data[17,4] = 0.01666666666666667
$data[18,0] = 1.701492537313433
$data[18,1] = 0.3766666666666666
$data[18,2] = 1
$data[18,3] = 0.2066666666666667
$data[18,4] = 0.01666666666666667
$data[19,0] = 1.736842105263158
$data[19,1] = 0.36
$data[19,2] = 1
$data[19,3] = 0.2366666666666667
$data[19,4] = 0.01666666666666667
$data[20,0] = 1.963855421686747
$data[20,1] = 0.3733333333333334
$data[20,2] = 1
$data[20,3] = 0.2566666666666667
$data[20,4] = 0.02
$data[21,0] = 1.943820224719101
$data[21,1] = 0.37
$data[21,2] = 1
$data[21,3] = 0.2766666666666667
$data[21,4] = 0.02
$data[22,0] = 1.928571428571429
$data[22,1] = 0.3866666666666667
$data[22,2] = 1
$data[22,3] = 0.3066666666666666
$data[22,4] = 0.02
$data[23,0] = 1.888888888888889
$data[23,1] = 0.3833333333333334
$data[23,2] = 1
$data[23,3] = 0.34
$data[23,4] = 0.02
$data[24,0] = 1.905982905982906
$data[24,1] = 0.3866666666666667
$data[24,2] = 1
$data[24,3] = 0.37
$data[24,4] = 0.02
$data[25,0] = 1.811023622047244
$data[25,1] = 0.3933333333333333
$data[25,2] = 1
$data[25,3] = 0.4033333333333333
$data[25,4] = 0.02
$data[26,0] = 1.876811594202898
$data[26,1] = 0.4066666666666667
$data[26,2] = 1
$data[26,3] = 0.4366666666666666
$data[26,4] = 0.02333333333333333
$data[27,0] = 1.888888888888889
$data[27,1] = 0.4066666666666667
$data[27,2] = 1
$data[27,3] = 0.4566666666666667
$data[27,4] = 0.02333333333333333
$data[28,0] = 1.855263157894737
$data[28,1] = 0.4
$data[28,2] = 1
$data[28,3] = 0.48
$data[28,4] = 0.02666666666666667
$data[29,0] = 1.895705521472393
$data[29,1] = 0.3666666666666666
$data[29,2] = 1
$data[29,3] = 0.5166666666666667
$data[29,4] = 0.02666666666666667
$data[30,0] = 1.861271676300578
$data[30,1] = 0.3333333333333333
$data[30,2] = 1
$data[30,3] = 0.55
$data[30,4] = 0.02666666666666667
$data[31,0] = 1.862637362637363
$data[31,1] = 0.31
$data[31,2] = 1
$data[31,3] = 0.58
$data[31,4] = 0.02666666666666667
$data[32,0] = 1.821989528795811
$data[32,1] = 0.2833333333333333
$data[32,2] = 1
$data[32,3] = 0.6066666666666667
$data[32,4] = 0.03
$data[33,0] = 1.84
$data[33,1] = 0.27
$data[33,2] = 1
$data[33,3] = 0.6333333333333333
$data[33,4] = 0.03333333333333333
$data[34,0] = 1.831730769230769
$data[34,1] = 0.25
$data[34,2] = 1
$data[34,3] = 0.66
$data[34,4] = 0.03333333333333333
$data[35,0] = 1.81042654028436
$data[35,1] = 0.2466666666666667
$data[35,2] = 1
$data[35,3] = 0.67
$data[35,4] = 0.03333333333333333
$data[36,0] = 1.784403669724771
$data[36,1] = 0.2333333333333333
$data[36,2] = 1
$data[36,3] = 0.6933333333333334
$data[36,4] = 0.03333333333333333
$data[37,0] = 1.836283185840708
$data[37,1] = 0.21
$data[37,2] = 1
$data[37,3] = 0.72
$data[37,4] = 0.03333333333333333
$data[38,0] = 1.807692307692308
$data[38,1] = 0.1866666666666667
$data[38,2] = 1
$data[38,3] = 0.7466666666666667
$data[38,4] = 0.03333333333333333
$data[39,0] = 1.79253112033195
$data[39,1] = 0.17
$data[39,2] = 1
$data[39,3] = 0.77
$data[39,4] = 0.03333333333333333
$data[40,0] = 1.766129032258065
$data[40,1] = 0.1566666666666667
$data[40,2] = 1
$data[40,3] = 0.7933333333333333
$data[40,4] = 0.03333333333333333
$data[41,0] = 1.782608695652174
$data[41,1] = 0.14
$data[41,2] = 1
$data[41,3] = 0.8100000000000001
$data[41,4] = 0.03333333333333333
$data[42,0] = 1.768339768339768
$data[42,1] = 0.12
$data[42,2] = 1
$data[42,3] = 0.83
$data[42,4] = 0.03333333333333333
$data[43,0] = 1.756653992395437
$data[43,1] = 0.1066666666666667
$data[43,2] = 1
$data[43,3] = 0.8433333333333334
$data[43,4] = 0.03333333333333333
$data[44,0] = 1.738805970149254
$data[44,1] = 0.09
$data[44,2] = 1
$data[44,3] = 0.86
$data[44,4] = 0.03333333333333333
$data[45,0] = 1.730627306273063
$data[45,1] = 0.08
$data[45,2] = 1
$data[45,3] = 0.87
$data[45,4] = 0.03333333333333333
$data[46,0] = 1.732600732600733
$data[46,1] = 0.07333333333333333
$data[46,2] = 1
$data[46,3] = 0.8766666666666667
$data[46,4] = 0.03333333333333333
$data[47,0] = 1.723636363636364
$data[47,1] = 0.06666666666666667
$data[47,2] = 1
$data[47,3] = 0.8833333333333333
$data[47,4] = 0.03333333333333333
$data[48,0] = 1.714801444043321
$data[48,1] = 0.06
$data[48,2] = 1
$data[48,3] = 0.8866666666666667
$data[48,4] = 0.03666666666666667
$data[49,0] = 1.707142857142857
$data[49,1] = 0.05
$data[49,2] = 1
$data[49,3] = 0.8966666666666666
$data[49,4] = 0.03666666666666667
$data[50,0] = 1.707142857142857
$data[50,1] = 0.05
$data[50,2] = 1
$data[50,3] = 0.8966666666666666
$data[50,4] = 0.03666666666666667
$data[51,0] = 1.692579505300353
$data[51,1] = 0.04
$data[51,2] = 1
$data[51,3] = 0.9066666666666666
$data[51,4] = 0.03666666666666667
$data[52,0] = 1.692579505300353
$data[52,1] = 0.04
$data[52,2] = 1
$data[52,3] = 0.9066666666666666
$data[52,4] = 0.03666666666666667
$data[53,0] = 1.687719298245614
$data[53,1] = 0.03333333333333333
$data[53,2] = 1
$data[53,3] = 0.9133333333333333
$data[53,4] = 0.03666666666666667
$data[54,0] = 1.70383275261324
$data[54,1] = 0.02666666666666667
$data[54,2] = 1
$data[54,3] = 0.92
$data[54,4] = 0.03666666666666667
$data[55,0] = 1.71875
$data[55,1] = 0.02333333333333333
$data[55,2] = 1
$data[55,3] = 0.9233333333333333
$data[55,4] = 0.03666666666666667
$data[56,0] = 1.717241379310345
$data[56,1] = 0.01666666666666667
$data[56,2] = 1
$data[56,3] = 0.93
$data[56,4] = 0.03666666666666667
$data[57,0] = 1.717241379310345
$data[57,1] = 0.01666666666666667
$data[57,2] = 1
$data[57,3] = 0.93
$data[57,4] = 0.03666666666666667
$data[58,0] = 1.717241379310345
$data[58,1] = 0.01666666666666667
$data[58,2] = 1
$data[58,3] = 0.93
$data[58,4] = 0.03666666666666667
$data[59,0] = 1.711340206185567
$data[59,1] = 0.01333333333333333
$data[59,2] = 1
$data[59,3] = 0.9333333333333333
$data[59,4] = 0.03666666666666667
$data[60,0] = 1.711340206185567
$data[60,1] = 0.01333333333333333
$data[60,2] = 1
$data[60,3] = 0.9333333333333333
$data[60,4] = 0.03666666666666667
$data[61,0] = 1.705479452054794
$data[61,1] = 0.01
$data[61,2] = 1
$data[61,3] = 0.9366666666666666
$data[61,4] = 0.03666666666666667
$data[62,0] = 1.704081632653061
$data[62,1] = 0.003333333333333334
$data[62,2] = 1
$data[62,3] = 0.9433333333333334
$data[62,4] = 0.03666666666666667
$data[63,0] = 1.704081632653061
$data[63,1] = 0.003333333333333334
$data[63,2] = 1
$data[63,3] = 0.9433333333333334
$data[63,4] = 0.03666666666666667
$data[64,0] = 0
$data[64,1] = 0
$data[64,2] = 0
$data[64,3] = 0
$data[64,4] = 0

$range = $ws.Range("B2:F66")
$range.Value = $data

Write-Host "Updated B2:F66 with new simulation values"